$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-26 Sunday" "2025-01-27 Monday"

Replace-Text "15×39=585" "76×96=7296"
Replace-Text "70×27=1890" "61×25=1525"
Replace-Text "25×90=2250" "21×17=357"
Replace-Text "37×60=2220" "36×99=3564"
Replace-Text "44×98=4312" "53×67=3551"

Replace-Text "38×11=418" "49×65=3185"
Replace-Text "46×97=4462" "32×86=2752"
Replace-Text "59×24=1416" "21×32=672"
Replace-Text "64×85=5440" "35×39=1365"
Replace-Text "74×43=3182" "80×67=5360"

Replace-Text "76×73=5548" "40×72=2880"
Replace-Text "91×32=2912" "81×30=2430"
Replace-Text "70×82=5740" "86×98=8428"
Replace-Text "54×76=4104" "34×66=2244"
Replace-Text "12×41=492" "14×71=994"

Replace-Text "53×69=3657" "83×93=7719"
Replace-Text "79×21=1659" "25×71=1775"
Replace-Text "56×94=5264" "76×90=6840"
Replace-Text "88×82=7216" "21×23=483"
Replace-Text "14×85=1190" "85×85=7225"

Replace-Text "89×68=6052" "38×59=2242"
Replace-Text "99×52=5148" "24×59=1416"
Replace-Text "66×83=5478" "82×67=5494"
Replace-Text "33×19=627" "65×46=2990"
Replace-Text "13×51=663" "57×21=1197"
